# Update "想去人数" (want-to-go count) figures pulled from the latest
# bilibili show-info scrape. Two sheets carry the same events: "展览"
# (exhibitions) and "全部类型" (all types combined).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7190
$ws1.Range("F4").Value = 116
$ws1.Range("F5").Value = 168
$ws1.Range("F7").Value = 87
$ws1.Range("F8").Value = 599
$ws1.Range("F9").Value = 52

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7190
$ws4.Range("F5").Value = 116
$ws4.Range("F6").Value = 168
$ws4.Range("F9").Value = 87
$ws4.Range("F10").Value = 599
$ws4.Range("F11").Value = 52
